# "table border is fixed"
# The data table (header + rows) on Sheet1 had no visible borders even though
# several cell styles referenced a border definition - that border definition
# was empty (no line on any side). Fix it by applying a thin border to every
# edge (outside and inside) of the table range, matching the sheet's used
# dimension A1:E14.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tableRange = $ws.Range("A1:E14")

# xlContinuous = 1, xlThin = 2
$tableRange.Borders.LineStyle = 1
$tableRange.Borders.Weight = 2
